$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-08 Saturday" "2025-02-09 Sunday"

Replace-Text "18÷2=9, 0" "80÷7=11, 3"
Replace-Text "78÷9=8, 6" "45÷4=11, 1"
Replace-Text "98÷7=14, 0" "62÷8=7, 6"
Replace-Text "84÷7=12, 0" "47÷6=7, 5"
Replace-Text "70÷6=11, 4" "89÷6=14, 5"

Replace-Text "30÷2=15, 0" "24÷5=4, 4"
Replace-Text "18÷5=3, 3" "65÷3=21, 2"
Replace-Text "40÷5=8, 0" "29÷6=4, 5"
Replace-Text "41÷6=6, 5" "59÷6=9, 5"
Replace-Text "38÷3=12, 2" "64÷7=9, 1"

Replace-Text "23÷9=2, 5" "29÷6=4, 5"
Replace-Text "59÷9=6, 5" "30÷4=7, 2"
Replace-Text "29÷9=3, 2" "16÷7=2, 2"
Replace-Text "77÷7=11, 0" "33÷7=4, 5"

Replace-Text "51÷7=7, 2" "14÷5=2, 4"
Replace-Text "56÷7=8, 0" "25÷3=8, 1"
Replace-Text "54÷8=6, 6" "39÷8=4, 7"
Replace-Text "76÷8=9, 4" "72÷3=24, 0"
Replace-Text "97÷6=16, 1" "10÷6=1, 4"

# This must run before the "10÷5=2, 0" -> "33÷6=5, 3" replacement below,
# since that replacement would otherwise create a duplicate "33÷6=5, 3"
# text earlier in the document and be matched here instead.
Replace-Text "33÷6=5, 3" "40÷4=10, 0"

Replace-Text "10÷5=2, 0" "33÷6=5, 3"

Replace-Text "13÷2=6, 1" "15÷3=5, 0"
Replace-Text "12÷8=1, 4" "82÷3=27, 1"
Replace-Text "65÷5=13, 0" "94÷6=15, 4"
Replace-Text "56÷6=9, 2" "36÷7=5, 1"
